# PlayerPerformance_6795.xlsx update:
#  1. Insert a new "Player Info" sheet before the existing "ODI Batting" sheet,
#     with player ID/NAME/BATTING_HAND/BOWL_STYLE data.
#  2. On the "ODI Batting" sheet, rename the MATCH_CARD_LINK column to
#     MATCH_CODE and replace the full scorecard URL with just the match code.

$wb = $excel.ActiveWorkbook

# --- 1. New "Player Info" sheet, inserted before "ODI Batting" -------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Look the batting sheet back up by name (its positional index shifted once
# the new sheet was inserted in front of it).
$odiBatting = $wb.Worksheets.Item("ODI Batting")

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $playerInfo.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "6795"
$playerInfo.Cells.Item(2, 1).Style = "Normal"
$playerInfo.Cells.Item(2, 2).Value = "Dane Cleaver"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Leg Break"

$playerInfo.PageSetup.LeftMargin = 0.75 * 72
$playerInfo.PageSetup.RightMargin = 0.75 * 72
$playerInfo.PageSetup.TopMargin = 1 * 72
$playerInfo.PageSetup.BottomMargin = 1 * 72
$playerInfo.PageSetup.HeaderMargin = 0.5 * 72
$playerInfo.PageSetup.FooterMargin = 0.5 * 72

# --- 2. Update "ODI Batting" sheet column D ---------------------------------
$odiBatting.Range("D1").Value = "MATCH_CODE"

$odiBatting.Range("D2").NumberFormat = "@"
$odiBatting.Range("D2").Value = "4625"
$odiBatting.Range("D2").Style = "Normal"
